# Add coefficient tables from multistate NMA and MA
# Adds d_name / mu_name lookup columns (G, H) to each of the 4 parameter
# lookup sheets (weibull, gompertz, fracpoly1, fracpoly2).

$wb = $excel.ActiveWorkbook

$wsWeibull   = $wb.Worksheets.Item("weibull")
$wsGompertz  = $wb.Worksheets.Item("gompertz")
$wsFracpoly1 = $wb.Worksheets.Item("fracpoly1")
$wsFracpoly2 = $wb.Worksheets.Item("fracpoly2")

# ---------------------------------------------------------------------
# weibull (sheet1): transition_id 3/4/5, mu_num sequence 1,2,3,(-),4,5
# ---------------------------------------------------------------------
$wsWeibull.Range("G1").Value = "d_name"
$wsWeibull.Range("H1").Value = "mu_name"
$wsWeibull.Range("H2").Value = "mu_1"
$wsWeibull.Range("H3").Value = "mu_2"
$wsWeibull.Range("H4").Value = "mu_4"
$wsWeibull.Range("H6").Value = "mu_5"
$wsWeibull.Range("H7").Value = "mu_6"

# ---------------------------------------------------------------------
# gompertz (sheet2): same row/transition layout as weibull
# ---------------------------------------------------------------------
$wsGompertz.Range("G1").Value = "d_name"
$wsGompertz.Range("H1").Value = "mu_name"
$wsGompertz.Range("H2").Value = "mu_1"
$wsGompertz.Range("H3").Value = "mu_2"
$wsGompertz.Range("H4").Value = "mu_4"
$wsGompertz.Range("H6").Value = "mu_5"
$wsGompertz.Range("H7").Value = "mu_6"

# ---------------------------------------------------------------------
# fracpoly1 (sheet3): transition_id 3/4/5 x gamma1/gamma2/gamma3,
# mu_num sequence 1,2,3,4,(-),(-),5,6,(-)
# ---------------------------------------------------------------------
$wsFracpoly1.Range("G1").Value = "d_name"
$wsFracpoly1.Range("H1").Value = "mu_name"
$wsFracpoly1.Range("H2").Value = "mu_1"
$wsFracpoly1.Range("H3").Value = "mu_2"
$wsFracpoly1.Range("H4").Value = "mu_3"
$wsFracpoly1.Range("H5").Value = "mu_4"
$wsFracpoly1.Range("H8").Value = "mu_5"
$wsFracpoly1.Range("H9").Value = "mu_6"

# ---------------------------------------------------------------------
# fracpoly2 (sheet4): identical layout to fracpoly1
# ---------------------------------------------------------------------
$wsFracpoly2.Range("G1").Value = "d_name"
$wsFracpoly2.Range("H1").Value = "mu_name"
$wsFracpoly2.Range("H2").Value = "mu_1"
$wsFracpoly2.Range("H3").Value = "mu_2"
$wsFracpoly2.Range("H4").Value = "mu_3"
$wsFracpoly2.Range("H5").Value = "mu_4"
$wsFracpoly2.Range("H8").Value = "mu_5"
$wsFracpoly2.Range("H9").Value = "mu_6"

# ---------------------------------------------------------------------
# View/selection state to match the authored workbook
# ---------------------------------------------------------------------
$wsWeibull.Range("G1:H1").Select()
$wsGompertz.Range("H2:H7").Select()
$wsFracpoly1.Range("H2:H9").Select()
$wsFracpoly2.Range("H2:H9").Select()

# fracpoly2 ends up the active/selected tab
$wsFracpoly2.Select()
